$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column before the existing "nom" column (BW),
# which shifts "nom" -> BX and "url_produit" -> BY.
$ws.Columns("BW:BW").Insert()

# Header for the freshly inserted price-check timestamp column.
$ws.Range("BW1").Value = "2026-01-31 03:13:36"

# Populate the new BW column for each product row by carrying forward
# the most recent known price from column BV (same value means price
# unchanged at this check); rows without a previous price stay blank.
for ($r = 2; $r -le 206; $r++) {
    $ws.Cells.Item($r, 75).Value = $ws.Cells.Item($r, 74).Value2
}
